# Adding 5 search test cases(B68-B72)
# (as captured by the canonical diff: two new rows appended to the
# "Test Cases" sheet, carrying six new shared strings.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 25 picks up the formatting already used by the last data row (24):
# thin border on every cell, no fill, no wrap - same visual style as the
# rest of the table.
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A25:E26").ClearContents()

$ws.Range("A25:E25").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)

# Row 25: TestCase_A24
# Values are entered A, then C, then B to match the shared-string
# insertion order recorded in the workbook.
$ws.Range("A25").Value = "TestCase_A24"
$ws.Range("C25").Value = "Verify that TR account gets locked after 5 consecutive unsuccessful login attempts"
$ws.Range("B25").Value = "OPQA-525"
$ws.Range("D25").Value = "Y"
$ws.Range("E25").Value = "SKIP"

# Row 26: TestCase_A25
$ws.Range("A26").Value = "TestCase_A25"
$ws.Range("B26").Value = "OPQA-529"
$ws.Range("C26").Value = "Verify that Help link is working properly"
$ws.Range("D26").Value = "Y"
$ws.Range("E26").Value = "PASS"

$ws.Range("D19").Select()
